$d = $word.ActiveDocument

function Get-ParaByPrefix($prefix) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs($i)
        if ($p.Range.Text.StartsWith($prefix)) {
            return $i
        }
    }
    return -1
}

function Rebuild-Paragraph($idx, $segments) {
    # $segments is an array of hashtables: @{ Text = "..."; Underline = $true/$false }
    # Strategy: clear the paragraph's content, insert the full concatenated text as a
    # single block (so Word doesn't propagate formatting from a previous insertion
    # into the next one), then apply Underline formatting to the specific sub-ranges
    # that need it, computed from absolute character offsets.
    $p = $d.Paragraphs($idx)
    $r = $p.Range
    $start = $r.Start
    $end = $r.End
    if (($end - 1) -gt $start) {
        $clearRange = $d.Range($start, $end - 1)
        $clearRange.Text = ""
    }

    $fullText = ""
    foreach ($seg in $segments) {
        $fullText = $fullText + $seg.Text
    }

    $ins = $d.Range($start, $start)
    $ins.InsertBefore($fullText)

    $pos = $start
    foreach ($seg in $segments) {
        $len = $seg.Text.Length
        if ($seg.Underline) {
            $fmtRange = $d.Range($pos, $pos + $len)
            $fmtRange.Font.Underline = 1
        }
        $pos = $pos + $len
    }
}

# --- Jogador ---
$idx = Get-ParaByPrefix("Jogador (")
Rebuild-Paragraph $idx @(
    @{ Text = "Jogador "; Underline = $false },
    @{ Text = "= "; Underline = $false },
    @{ Text = "{"; Underline = $false },
    @{ Text = "Nome"; Underline = $true },
    @{ Text = ", Sexo, Raça, Classe"; Underline = $false },
    @{ Text = ", "; Underline = $false },
    @{ Text = "Nível, Força de Combate"; Underline = $false },
    @{ Text = ", "; Underline = $false },
    @{ Text = "Equipados: [Cabeça"; Underline = $false },
    @{ Text = ", Mão Direita, Mão Esquerda, Duas Mãos, Armadura, Pés"; Underline = $false },
    @{ Text = "]"; Underline = $false },
    @{ Text = "}"; Underline = $false }
)

# --- Mochila ---
$idx = Get-ParaByPrefix("Mochila (")
Rebuild-Paragraph $idx @(
    @{ Text = "Mochila"; Underline = $false },
    @{ Text = " ="; Underline = $false },
    @{ Text = " "; Underline = $false },
    @{ Text = "{"; Underline = $false },
    @{ Text = "id"; Underline = $true },
    @{ Text = ", Capacidade"; Underline = $false },
    @{ Text = "}"; Underline = $false }
)

# --- Item ---
$idx = Get-ParaByPrefix("Item (")
Rebuild-Paragraph $idx @(
    @{ Text = "Item "; Underline = $false },
    @{ Text = "= "; Underline = $false },
    @{ Text = "{"; Underline = $false },
    @{ Text = "Nome"; Underline = $true },
    @{ Text = ", Valor, Bônus"; Underline = $false },
    @{ Text = "}"; Underline = $false }
)

# --- Equipamento ---
$idx = Get-ParaByPrefix("Equipamento (")
Rebuild-Paragraph $idx @(
    @{ Text = "Equipamento"; Underline = $false },
    @{ Text = " ="; Underline = $false },
    @{ Text = " "; Underline = $false },
    @{ Text = "{"; Underline = $false },
    @{ Text = "N"; Underline = $true },
    @{ Text = "ome"; Underline = $true },
    @{ Text = ", Valor, Bônus, "; Underline = $false },
    @{ Text = "Tamanho"; Underline = $false },
    @{ Text = ", "; Underline = $false },
    @{ Text = "Local Ocupado"; Underline = $false },
    @{ Text = "}"; Underline = $false }
)

# --- Consumível ---
$idx = Get-ParaByPrefix("Consumível (")
Rebuild-Paragraph $idx @(
    @{ Text = "Consumível "; Underline = $false },
    @{ Text = "= "; Underline = $false },
    @{ Text = "{"; Underline = $false },
    @{ Text = "N"; Underline = $true },
    @{ Text = "ome"; Underline = $true },
    @{ Text = ", Valor, Bônus, "; Underline = $false },
    @{ Text = "Efeito"; Underline = $false },
    @{ Text = "}"; Underline = $false }
)

# --- Monstro ---
$idx = Get-ParaByPrefix("Monstro (")
Rebuild-Paragraph $idx @(
    @{ Text = "Monstro "; Underline = $false },
    @{ Text = "= "; Underline = $false },
    @{ Text = "{"; Underline = $false },
    @{ Text = "Nome"; Underline = $true },
    @{ Text = ", "; Underline = $false },
    @{ Text = "Nível,"; Underline = $false },
    @{ Text = " D"; Underline = $false },
    @{ Text = "rops, Quantidade de níveis ganhos ao derrotar"; Underline = $false },
    @{ Text = "}"; Underline = $false }
)

# --- Sala ---
$idx = Get-ParaByPrefix("Sala (")
if ($idx -eq -1) {
    $idx = Get-ParaByPrefix("Sala ")
}
Rebuild-Paragraph $idx @(
    @{ Text = "Sala "; Underline = $false },
    @{ Text = "= "; Underline = $false },
    @{ Text = "{"; Underline = $false },
    @{ Text = "N"; Underline = $false },
    @{ Text = "ome"; Underline = $false },
    @{ Text = "}"; Underline = $false }
)

Write-Host "Done"
